# Update the two-digit multiplication problems in the table.
# Each non-empty row of the 5-column table gets new operands; addressed
# by Table.Cell(row, col) so that duplicate text (e.g. "88×73=" which
# appears twice in the original) is updated independently and correctly.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("33×47=", "73×68=", "62×42=", "79×95=", "62×47=")
    5  = @("72×55=", "44×28=", "70×23=", "43×16=", "42×97=")
    10 = @("32×96=", "19×26=", "27×83=", "19×56=", "34×41=")
    15 = @("52×20=", "58×32=", "40×39=", "66×13=", "29×92=")
    20 = @("78×98=", "40×88=", "39×97=", "74×57=", "78×28=")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($col = 1; $col -le $values.Length; $col++) {
        $t.Cell($row, $col).Range.Text = $values[$col - 1]
    }
}
